$d = $word.ActiveDocument

# The document ends with: </w:tbl><w:p/><w:p/><w:sectPr>...
# We keep the first trailing empty paragraph as-is, and turn the second
# (last) paragraph into a new "Reviderad tidsplan" heading, followed by
# two new body paragraphs, inserted just before the final section break.

$n = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($n)
$headingPara.Range.Text = "Reviderad tidsplan"
$headingPara.Range.Style = "Rubrik1"

$headingPara.Range.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p2 = $d.Paragraphs.Item($n)
$p2.Range.Style = "Normal"
$p2.Range.Text = "Efter att ha arbetat en del med projektet så vet vi att planeringen inte riktigt håller. Grafiken har tagit mycket mer tid än vi förväntat och vi har fortfarande inte gjort kollision eller börjat med nätverk."

$p2.Range.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p3 = $d.Paragraphs.Item($n)
$p3.Range.Style = "Normal"
$p3.Range.Text = "Vi kommer troligtvis inte hinna att färdigställa spelet, eller få mycket mer än en demo för själva spelmotorn, men vi kommer förhoppningsvis hinna med att färdigställa kollisionen samt nätverks delen. Däremot får vi börja med projektrapporten snart, helst en bit innan sista veckan, så att vi då endast behöver uppdatera det nya vi gjort."

Write-Output "paragraph count: $($d.Paragraphs.Count)"
Write-Output "Content end: $($d.Content.End)"
